$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '63.627.73'
$ws.Range('E2').Value = '  +0.08%  '

# Row 3
$ws.Range('D3').Value = '2.651.65'
$ws.Range('E3').Value = '  +0.02%  '

# Row 4
$ws.Range('E4').Value = '  +0.10%  '

# Row 5
$ws.Range('D5').Value = '''602.05'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.81%  '

# Row 6
$ws.Range('D6').Value = '''147.40'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.94%  '

# Row 7
$ws.Range('E7').Value = '  +0.10%  '

# Row 8
$ws.Range('D8').Value = '''0.588'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.02%  '

# Row 9
$ws.Range('D9').Value = '''0.108'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.57%  '

# Row 10
$ws.Range('D10').Value = '''5.59'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.37%  '

# Row 11
$ws.Range('D11').Value = '''0.370'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +4.59%  '

# Row 12
$ws.Range('E12').Value = '  +0.04%  '

# Row 13
$ws.Range('D13').Value = '''27.53'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.29%  '

# Row 14
$ws.Range('D14').Value = '3.129.46'
$ws.Range('E14').Value = '  +0.18%  '

# Row 15
$ws.Range('D15').Value = '63.507.61'
$ws.Range('E15').Value = '  +0.04%  '

# Row 16
$ws.Range('E16').Value = '  +0.34%  '

# Row 17
$ws.Range('D17').Value = '2.646.88'
$ws.Range('E17').Value = '  +0.54%  '

# Row 18
$ws.Range('D18').Value = '''11.46'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.51%  '

# Row 19
$ws.Range('D19').Value = '''4.55'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +4.45%  '

# Row 20
$ws.Range('D20').Value = '''342.44'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.55%  '

# Row 21
$ws.Range('D21').Value = '''6.96'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +3.44%  '

# Row 22
$ws.Range('E22').Value = '  -0.02%  '

# Row 23
$ws.Range('D23').Value = '''5.57'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -3.35%  '

# Row 24
$ws.Range('D24').Value = '''66.85'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -1.11%  '

# Row 25
$ws.Range('D25').Value = '''1.69'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +3.11%  '

# Row 26
$ws.Range('D26').Value = '''9.05'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +7.20%  '

# Row 27
$ws.Range('D27').Value = '''1.54'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.89%  '

# Row 28
$ws.Range('D28').Value = '''558.62'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.37%  '

# Row 29
$ws.Range('D29').Value = '''0.164'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.87%  '

# Row 30
$ws.Range('D30').Value = '''1.00'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +0.13%  '

# Row 31
$ws.Range('D31').Value = '''7.88'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.34%  '

# Row 32
$ws.Range('D32').Value = '''2.03'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +2.48%  '

# Row 33
$ws.Range('D33').Value = '''1.75'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -3.84%  '

# Row 34
$ws.Range('D34').Value = '0.0₃0819'
$ws.Range('E34').Value = '  +1.16%  '

# Row 35
$ws.Range('D35').Value = '''5.19'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +6.97%  '

# Row 36
$ws.Range('D36').Value = '''166.63'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -4.46%  '

# Row 37
$ws.Range('D37').Value = '''0.406'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.87%  '

# Row 38
$ws.Range('E38').Value = '  +0.09%  '

# Row 39
$ws.Range('E39').Value = '  +6.75%  '

# Row 40
$ws.Range('D40').Value = '''19.13'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.26%  '

# Row 41
$ws.Range('E41').Value = '  +0.02%  '

# Row 42
$ws.Range('D42').Value = '''167.84'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.53%  '

# Row 43
$ws.Range('D43').Value = '''3.79'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +1.29%  '

# Row 44
$ws.Range('D44').Value = '''22.24'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.08%  '

# Row 45
$ws.Range('D45').Value = '''0.0575'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +3.98%  '

# Row 46
$ws.Range('D46').Value = '''0.629'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.24%  '

# Row 47
$ws.Range('D47').Value = '''0.0248'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +3.76%  '

# Row 48
$ws.Range('D48').Value = '''0.0962'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.14%  '

# Row 49
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '''18.80'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.33%  '

# Row 50
$ws.Range('B50').Value = 'dogwifhat'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D50').Value = '''1.87'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +9.87%  '

# Row 51
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').Value = '''11.28'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.50%  '

